$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sample" column (A) for the spleen rows of this FACS sheet was missing
# the "<n>: " index prefix that every other sample group already has (e.g.
# "1: mLN_681_030.fcs"). Add the missing "1: " prefix to each of those
# spleen_*.fcs sample names (rows 147-191) so the sheet is consistent and
# can be parsed programmatically (e.g. in R).
for ($row = 147; $row -le 191; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = [string]$cell.Text
    if ($current.Length -gt 0 -and -not ($current -match '^\d+:\s')) {
        $cell.Value = "1: " + $current
    }
}

# Restore the sheet's active cell/view state to where the editor left off.
$ws.Range("C154").Select()
